# İş Takip Güncellemesi - 24.07.2025 15:43:15
# Populate the "İş Takip Listesi" sheet (header row + first data row)
# and keep every value stored as text, matching the source export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("A1:L2")
# Force text formatting first so numeric/date-looking strings (parcel
# counts, areas, dates) are not silently converted to numbers/dates.
$range.NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "Tamamlandı"
$ws.Range("B1").Value = "İL"
$ws.Range("C1").Value = "İLÇE"
$ws.Range("D1").Value = "BİRİM"
$ws.Range("E1").Value = "UYGULAMA"
$ws.Range("F1").Value = "GÖREVLİ PERSONELLER"
$ws.Range("G1").Value = "PARSEL SAYISI"
$ws.Range("H1").Value = "ALAN(Ha)"
$ws.Range("I1").Value = "İHALELİ/MÜDÜRLÜK"
$ws.Range("J1").Value = "İŞE BAŞLAMA/YER TESLİMİ"
$ws.Range("K1").Value = "İHALE BİTİŞ TARİHİ"
$ws.Range("L1").Value = "DURUMU"

# Data row
$ws.Range("A2").Value = "HAYIR"
$ws.Range("B2").Value = "Adana"
$ws.Range("C2").Value = "Akdeniz"
$ws.Range("D2").Value = "gggg"
$ws.Range("E2").Value = "GÜNCELLEME"
$ws.Range("F2").Value = "EMİNE ALANLI KIRCILI (K.Mühendisi), HİLMİ MÜFTÜOĞLU (K.Mühendisi)"
$ws.Range("G2").Value = "3"
$ws.Range("H2").Value = "120"
$ws.Range("I2").Value = "İhaleli"
$ws.Range("J2").Value = "2025-02-01"
$ws.Range("K2").Value = "2025-07-24"
$ws.Range("L2").Value = "ARAZİ DEVAM EDİYOR"
